$d = $word.ActiveDocument

# Locate the target paragraph: the empty paragraph carrying the stray
# <w:ind w:left="708" w:hanging="708"/> (708 twips = 35.4 pt) left indent,
# which sits right after the Sprint 2 bullet list ("Sprint 3" section goes here).
$target = $null
$matches = 0
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13)
    $li = $p.LeftIndent
    if ($txt -eq "" -and [Math]::Abs($li - 35.4) -lt 0.1) {
        $target = $p
        $matches = $matches + 1
    }
}

if ($target -eq $null -or $matches -ne 1) {
    throw "Could not uniquely locate target paragraph (matches=$matches)"
}

# Replace that paragraph (pilcrow included) with the new Sprint 3 header,
# its two bullet requirements, and a fresh trailing blank paragraph.
# (Two trailing <w:p/> are needed in the payload: InsertXML's own final
# paragraph mark always fuses with the destination's original mark, so an
# extra one is required to actually get a new blank paragraph in the result.)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">Sprint 3 - </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Creación de la base de datos y la estructura arquitectónica del proyecto</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> (1 Semana) 15 – 21 enero 2024</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Generar la división de capas (arquitectura) en el proyecto de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Netbeans</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Crear la estructura de la base de datos con el SGBD MySQL.</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.Range.InsertXML($xml)

Write-Output "Sprint 3 section inserted"
